$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Terminal La Palmera de La Serena -
# Perejil". It belongs right above the current row 66 (by date ordering used
# in the sheet), so insert a new row there; this pushes the former rows
# 66-97 down to 67-98 while keeping their data intact.
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new weekly record.
$ws.Range("A66").Value = 8
$ws.Range("B66").Value = "Terminal La Palmera de La Serena"
$ws.Range("C66").Value = "Coquimbo"
$ws.Range("D66").Value = 44523
$ws.Range("E66").Value = 4
$ws.Range("F66").Value = 100112044
$ws.Range("G66").Value = "Perejil"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 2800
$ws.Range("K66").Value = 1500
$ws.Range("L66").Value = 2000
$ws.Range("M66").Value = 1750
$ws.Range("N66").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O66").Value = "Provincia del Elquí"
$ws.Range("P66").Value = 1167
$ws.Range("Q66").Value = 1.5
$ws.Range("R66").Value = "Hortaliza"
